# Finalize the project with documents.
#
# This script brings the Tasks workbook up to date:
#   - Modeling sheet: fill in the "Class diagram" note and the remaining
#     "Performer" entries for the sequence diagrams.
#   - Coding sheet: tweak the "Login form" note (login -> login/logout) and
#     add a note to the "Checkout form" row; becomes the active tab.
#   - Documenting sheet: record who is doing the "STC conclusion" task.
#
# NL = newline used inside wrapped note cells (matches how Excel stores
# multi-line text inside a shared string).
$NL = [char]10

$wb = $excel.ActiveWorkbook

$wsModeling    = $wb.Worksheets.Item("Modeling")
$wsCoding      = $wb.Worksheets.Item("Coding")
$wsDocumenting = $wb.Worksheets.Item("Documenting")

# ---------------------------------------------------------------------
# Modeling sheet
# ---------------------------------------------------------------------

# Row 3 "Class diagram": add a Note. Use the same wrapped/quoted cell
# style already used for note cells elsewhere in the workbook (copy format
# only, so the shared cell style is reused instead of duplicated).
$wsModeling.Range("E3").Value = "- one general class diagram" + $NL + "- one for new member use case"
$wsCoding.Range("E2").Copy()
$wsModeling.Range("E3").PasteSpecial(-4122)
$wsModeling.Rows.Item(3).RowHeight = 28.5

# Rows 4-6 "Sequence diagram 1/2/3": fill in the Performer column.
$wsModeling.Range("D4").Value = "Tin"
$wsModeling.Range("D5").Value = "Tin"
$wsModeling.Range("D6").Value = "Tin"

# ---------------------------------------------------------------------
# Coding sheet
# ---------------------------------------------------------------------

# Row 3 "Login form": note now also mentions logout. Re-apply the
# wrapped/quoted note style afterwards since writing a multi-line value
# resets the cell to a plain wrap style.
$wsCoding.Range("E3").Value = "- login/logout" + $NL + "- show/hide corresponding menus"
$wsCoding.Range("E2").Copy()
$wsCoding.Range("E3").PasteSpecial(-4122)

# Row 6 "Checkout form": add a Note (previously empty), reusing the
# wrapped/quoted note style.
$wsCoding.Range("E6").Value = "- input validation" + $NL + "- check availability of book"
$wsCoding.Range("E2").Copy()
$wsCoding.Range("E6").PasteSpecial(-4122)
$wsCoding.Rows.Item(6).RowHeight = 28.5

# ---------------------------------------------------------------------
# Documenting sheet
# ---------------------------------------------------------------------

# Row 3 "STC conclusion" task: Assignee + Performer.
$wsDocumenting.Range("B3").Value = "STC conclusion"
$wsDocumenting.Range("D3").Value = "Tin"

# ---------------------------------------------------------------------
# Selection / active tab bookkeeping (matches the saved view state)
# ---------------------------------------------------------------------

$wsModeling.Range("B7").Select()
$wsDocumenting.Range("D4").Select()
$wsCoding.Activate()
$wsCoding.Range("E7").Select()

$excel.CutCopyMode = $false
